$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.05707143610673
$ws.Cells.Item(2, 3).Value = 10.11266727046377
$ws.Cells.Item(2, 4).Value = 4.866972458720801
$ws.Cells.Item(2, 6).Value = 22.91487345060163
$ws.Cells.Item(2, 7).Value = 26.36220539376198
$ws.Cells.Item(2, 8).Value = 13.76516096934731
$ws.Cells.Item(2, 9).Value = 21.74193069949498
$ws.Cells.Item(2, 11).Value = 7.925055994073841
$ws.Cells.Item(2, 12).Value = 10.67841275826011
$ws.Cells.Item(2, 14).Value = 18.19434022731312
$ws.Cells.Item(2, 15).Value = 20.61882427973013

$ws.Cells.Item(3, 2).Value = 10.76010400995305
$ws.Cells.Item(3, 3).Value = 10.13158757623668
$ws.Cells.Item(3, 4).Value = 4.784708140403266
$ws.Cells.Item(3, 6).Value = 22.94466011048424
$ws.Cells.Item(3, 7).Value = 26.41347593726986
$ws.Cells.Item(3, 8).Value = 13.80505345436734
$ws.Cells.Item(3, 9).Value = 21.82022024763959
$ws.Cells.Item(3, 11).Value = 7.645495042701882
$ws.Cells.Item(3, 12).Value = 10.64966426284792
$ws.Cells.Item(3, 14).Value = 18.24537806026862
$ws.Cells.Item(3, 15).Value = 20.68154893016588

$ws.Cells.Item(4, 2).Value = 10.57514588943297
$ws.Cells.Item(4, 3).Value = 10.14393968284393
$ws.Cells.Item(4, 4).Value = 4.732723468180642
$ws.Cells.Item(4, 6).Value = 22.96853066809147
$ws.Cells.Item(4, 7).Value = 26.4527335548402
$ws.Cells.Item(4, 8).Value = 13.83146104954275
$ws.Cells.Item(4, 9).Value = 21.87209013475012
$ws.Cells.Item(4, 11).Value = 7.513804815899531
$ws.Cells.Item(4, 12).Value = 10.63410903432772
$ws.Cells.Item(4, 14).Value = 18.27825821940866
$ws.Cells.Item(4, 15).Value = 20.72398550102336

$ws.Cells.Item(5, 2).Value = 10.4992296267838
$ws.Cells.Item(5, 3).Value = 10.14915859808014
$ws.Cells.Item(5, 4).Value = 4.711181993063521
$ws.Cells.Item(5, 6).Value = 22.97966025799769
$ws.Cells.Item(5, 7).Value = 26.47068156952797
$ws.Cells.Item(5, 8).Value = 13.84270376786594
$ws.Cells.Item(5, 9).Value = 21.89418268095699
$ws.Cells.Item(5, 11).Value = 7.459813876360849
$ws.Cells.Item(5, 12).Value = 10.62830229098408
$ws.Cells.Item(5, 14).Value = 18.29204610984758
$ws.Cells.Item(5, 15).Value = 20.74226425100983

$ws.Cells.Item(6, 2).Value = 10.48659437552132
$ws.Cells.Item(6, 3).Value = 10.15003640384787
$ws.Cells.Item(6, 4).Value = 4.707583869091176
$ws.Cells.Item(6, 6).Value = 22.98159295365835
$ws.Cells.Item(6, 7).Value = 26.47377943673955
$ws.Cells.Item(6, 8).Value = 13.84459969616033
$ws.Cells.Item(6, 9).Value = 21.89790881001531
$ws.Cells.Item(6, 11).Value = 7.450823000294863
$ws.Cells.Item(6, 12).Value = 10.62737036057343
$ws.Cells.Item(6, 14).Value = 18.29435910101871
$ws.Cells.Item(6, 15).Value = 20.74535891326379

$ws.Cells.Item(7, 2).Value = 10.57412410452124
$ws.Cells.Item(7, 3).Value = 10.14400931586668
$ws.Cells.Item(7, 4).Value = 4.732434380945092
$ws.Cells.Item(7, 6).Value = 22.96867509059825
$ws.Cells.Item(7, 7).Value = 26.45296771949509
$ws.Cells.Item(7, 8).Value = 13.83161072322594
$ws.Cells.Item(7, 9).Value = 21.87238421565015
$ws.Cells.Item(7, 11).Value = 7.513078451569136
$ws.Cells.Item(7, 12).Value = 10.63402856176444
$ws.Cells.Item(7, 14).Value = 18.2784425914739
$ws.Cells.Item(7, 15).Value = 20.72422802582208

$ws.Cells.Item(8, 2).Value = 10.95528714028517
$ws.Cells.Item(8, 3).Value = 10.11903876241484
$ws.Cells.Item(8, 4).Value = 4.838921495522266
$ws.Cells.Item(8, 6).Value = 22.92398481336423
$ws.Cells.Item(8, 7).Value = 26.37826626739361
$ws.Cells.Item(8, 8).Value = 13.77851891763362
$ws.Cells.Item(8, 9).Value = 21.76813608689932
$ws.Cells.Item(8, 11).Value = 7.830126195663444
$ws.Cells.Item(8, 12).Value = 10.6680682865807
$ws.Cells.Item(8, 14).Value = 18.21161850073349
$ws.Cells.Item(8, 15).Value = 20.63963681977354

$ws.Cells.Item(9, 2).Value = 11.67728781114774
$ws.Cells.Item(9, 3).Value = 10.07587986045792
$ws.Cells.Item(9, 4).Value = 5.035490723479407
$ws.Cells.Item(9, 6).Value = 22.88067188162898
$ws.Cells.Item(9, 7).Value = 26.29369399434883
$ws.Cells.Item(9, 8).Value = 13.68957684460572
$ws.Cells.Item(9, 9).Value = 21.59387042309674
$ws.Cells.Item(9, 11).Value = 8.48718788966351
$ws.Cells.Item(9, 12).Value = 10.75121423692075
$ws.Cells.Item(9, 14).Value = 18.09276994990978
$ws.Cells.Item(9, 15).Value = 20.50492970805129

$ws.Cells.Item(10, 2).Value = 12.18646045883521
$ws.Cells.Item(10, 3).Value = 10.04767979342024
$ws.Cells.Item(10, 4).Value = 5.1717423532733
$ws.Cells.Item(10, 6).Value = 22.87589111704933
$ws.Cells.Item(10, 7).Value = 26.26953493989521
$ws.Cells.Item(10, 8).Value = 13.6334636666386
$ws.Cells.Item(10, 9).Value = 21.48424664653889
$ws.Cells.Item(10, 11).Value = 8.932246915858038
$ws.Cells.Item(10, 12).Value = 10.8219424694787
$ws.Cells.Item(10, 14).Value = 18.01281788722626
$ws.Cells.Item(10, 15).Value = 20.42502577452772

$ws.Cells.Item(11, 2).Value = 12.4124056667384
$ws.Cells.Item(11, 3).Value = 10.03560597419769
$ws.Cells.Item(11, 4).Value = 5.231817967303675
$ws.Cells.Item(11, 6).Value = 22.87957968243465
$ws.Cells.Item(11, 7).Value = 26.26682026089282
$ws.Cells.Item(11, 8).Value = 13.60993774579822
$ws.Cells.Item(11, 9).Value = 21.43837776090368
$ws.Cells.Item(11, 11).Value = 9.12603981598266
$ws.Cells.Item(11, 12).Value = 10.85612680827613
$ws.Cells.Item(11, 14).Value = 17.97803077130699
$ws.Cells.Item(11, 15).Value = 20.39282675340315

$ws.Cells.Item(12, 2).Value = 12.49706941786832
$ws.Cells.Item(12, 3).Value = 10.03114191346351
$ws.Cells.Item(12, 4).Value = 5.254282464402438
$ws.Cells.Item(12, 6).Value = 22.88181775905892
$ws.Cells.Item(12, 7).Value = 26.26698302334145
$ws.Cells.Item(12, 8).Value = 13.60131649103186
$ws.Cells.Item(12, 9).Value = 21.42158397782757
$ws.Cells.Item(12, 11).Value = 9.19814376957326
$ws.Cells.Item(12, 12).Value = 10.86935261122946
$ws.Cells.Item(12, 14).Value = 17.96508445091074
$ws.Cells.Item(12, 15).Value = 20.38123138164483

$ws.Cells.Item(13, 2).Value = 12.47887664602084
$ws.Cells.Item(13, 3).Value = 10.03209853144598
$ws.Cells.Item(13, 4).Value = 5.249457165427892
$ws.Cells.Item(13, 6).Value = 22.88129836876907
$ws.Cells.Item(13, 7).Value = 26.26689500950489
$ws.Cells.Item(13, 8).Value = 13.6031604479333
$ws.Cells.Item(13, 9).Value = 21.42517520184504
$ws.Cells.Item(13, 11).Value = 9.182672381253001
$ws.Cells.Item(13, 12).Value = 10.8664918359533
$ws.Cells.Item(13, 14).Value = 17.96786259956741
$ws.Cells.Item(13, 15).Value = 20.38370205492889

$ws.Cells.Item(14, 2).Value = 12.4193894100055
$ws.Cells.Item(14, 3).Value = 10.03523655076691
$ws.Cells.Item(14, 4).Value = 5.233671902800212
$ws.Cells.Item(14, 6).Value = 22.87974696243554
$ws.Cells.Item(14, 7).Value = 26.2668097863579
$ws.Cells.Item(14, 8).Value = 13.60922270879353
$ws.Cells.Item(14, 9).Value = 21.43698458333801
$ws.Cells.Item(14, 11).Value = 9.131997712273582
$ws.Cells.Item(14, 12).Value = 10.85720932030831
$ws.Cells.Item(14, 14).Value = 17.97696112931179
$ws.Cells.Item(14, 15).Value = 20.39186081111162

$ws.Cells.Item(15, 2).Value = 12.38283272096142
$ws.Cells.Item(15, 3).Value = 10.03717273341908
$ws.Cells.Item(15, 4).Value = 5.223965571693159
$ws.Cells.Item(15, 6).Value = 22.87890617683812
$ws.Cells.Item(15, 7).Value = 26.26691266061255
$ws.Cells.Item(15, 8).Value = 13.61297345578301
$ws.Cells.Item(15, 9).Value = 21.44429317154751
$ws.Cells.Item(15, 11).Value = 9.100790207166501
$ws.Cells.Item(15, 12).Value = 10.85155985129954
$ws.Cells.Item(15, 14).Value = 17.98256375184651
$ws.Cells.Item(15, 15).Value = 20.39693615048177

$ws.Cells.Item(16, 2).Value = 12.17157310432233
$ws.Cells.Item(16, 3).Value = 10.04848398889051
$ws.Cells.Item(16, 4).Value = 5.167776983175957
$ws.Cells.Item(16, 6).Value = 22.87576790363551
$ws.Cells.Item(16, 7).Value = 26.26987896145549
$ws.Cells.Item(16, 8).Value = 13.63504137968125
$ws.Cells.Item(16, 9).Value = 21.48732482162982
$ws.Cells.Item(16, 11).Value = 8.919404599515071
$ws.Cells.Item(16, 12).Value = 10.81974820385329
$ws.Cells.Item(16, 14).Value = 18.0151230937539
$ws.Cells.Item(16, 15).Value = 20.42721364415893

$ws.Cells.Item(17, 2).Value = 12.04045811941072
$ws.Cells.Item(17, 3).Value = 10.05561599699423
$ws.Cells.Item(17, 4).Value = 5.132811080537921
$ws.Cells.Item(17, 6).Value = 22.87534331741737
$ws.Cells.Item(17, 7).Value = 26.27381908918064
$ws.Cells.Item(17, 8).Value = 13.64909155336974
$ws.Cells.Item(17, 9).Value = 21.51474817605835
$ws.Cells.Item(17, 11).Value = 8.805884231969015
$ws.Cells.Item(17, 12).Value = 10.80074192524293
$ws.Cells.Item(17, 14).Value = 18.03550212023594
$ws.Cells.Item(17, 15).Value = 20.44685146246549

$ws.Cells.Item(18, 2).Value = 11.96451393125782
$ws.Cells.Item(18, 3).Value = 10.05978918731997
$ws.Cells.Item(18, 4).Value = 5.112520823972004
$ws.Cells.Item(18, 6).Value = 22.87565114609717
$ws.Cells.Item(18, 7).Value = 26.27686428580586
$ws.Cells.Item(18, 8).Value = 13.65736112213342
$ws.Cells.Item(18, 9).Value = 21.5308977796887
$ws.Cells.Item(18, 11).Value = 8.739776332691978
$ws.Cells.Item(18, 12).Value = 10.78999968993762
$ws.Cells.Item(18, 14).Value = 18.04737269458215
$ws.Cells.Item(18, 15).Value = 20.45853715026889

$ws.Cells.Item(19, 2).Value = 11.93871201037655
$ws.Cells.Item(19, 3).Value = 10.06121437474892
$ws.Cells.Item(19, 4).Value = 5.105620522440601
$ws.Cells.Item(19, 6).Value = 22.87585023035156
$ws.Cells.Item(19, 7).Value = 26.2780290804598
$ws.Cells.Item(19, 8).Value = 13.66019339739092
$ws.Cells.Item(19, 9).Value = 21.5364303881038
$ws.Cells.Item(19, 11).Value = 8.717254675926224
$ws.Cells.Item(19, 12).Value = 10.78639536883064
$ws.Cells.Item(19, 14).Value = 18.05141750223099
$ws.Cells.Item(19, 15).Value = 20.46256076673283

$ws.Cells.Item(20, 2).Value = 12.0544710220784
$ws.Cells.Item(20, 3).Value = 10.05484943277375
$ws.Cells.Item(20, 4).Value = 5.136551845143518
$ws.Cells.Item(20, 6).Value = 22.87533139261326
$ws.Cells.Item(20, 7).Value = 26.27331902973994
$ws.Cells.Item(20, 8).Value = 13.64757640302448
$ws.Cells.Item(20, 9).Value = 21.5117899509109
$ws.Cells.Item(20, 11).Value = 8.818053181910503
$ws.Cells.Item(20, 12).Value = 10.80274560259769
$ws.Cells.Item(20, 14).Value = 18.03331731359646
$ws.Cells.Item(20, 15).Value = 20.44472055699445

$ws.Cells.Item(21, 2).Value = 12.43688719350198
$ws.Cells.Item(21, 3).Value = 10.03431191016258
$ws.Cells.Item(21, 4).Value = 5.238316229887296
$ws.Cells.Item(21, 6).Value = 22.88017983352914
$ws.Cells.Item(21, 7).Value = 26.26680250089507
$ws.Cells.Item(21, 8).Value = 13.60743427419248
$ws.Cells.Item(21, 9).Value = 21.43350025043672
$ws.Cells.Item(21, 11).Value = 9.146917117534059
$ws.Cells.Item(21, 12).Value = 10.85992826193078
$ws.Cells.Item(21, 14).Value = 17.97428252222718
$ws.Cells.Item(21, 15).Value = 20.38944815789659

$ws.Cells.Item(22, 2).Value = 12.68155903744822
$ws.Cells.Item(22, 3).Value = 10.02151893432681
$ws.Cells.Item(22, 4).Value = 5.303160379544719
$ws.Cells.Item(22, 6).Value = 22.88825096689378
$ws.Cells.Item(22, 7).Value = 26.26948380154713
$ws.Cells.Item(22, 8).Value = 13.58287478561365
$ws.Cells.Item(22, 9).Value = 21.38568971997874
$ws.Cells.Item(22, 11).Value = 9.354367884182386
$ws.Cells.Item(22, 12).Value = 10.89893391894922
$ws.Cells.Item(22, 14).Value = 17.9370214614255
$ws.Cells.Item(22, 15).Value = 20.35680868827064

$ws.Cells.Item(23, 2).Value = 12.55147885110255
$ws.Cells.Item(23, 3).Value = 10.02828934214834
$ws.Cells.Item(23, 4).Value = 5.26870754721405
$ws.Cells.Item(23, 6).Value = 22.88349546308703
$ws.Cells.Item(23, 7).Value = 26.26741772875116
$ws.Cells.Item(23, 8).Value = 13.59582936287186
$ws.Cells.Item(23, 9).Value = 21.41089978552199
$ws.Cells.Item(23, 11).Value = 9.244342337114391
$ws.Cells.Item(23, 12).Value = 10.87796917000366
$ws.Cells.Item(23, 14).Value = 17.95678775859814
$ws.Cells.Item(23, 15).Value = 20.37390987821426

$ws.Cells.Item(24, 2).Value = 12.04813754071553
$ws.Cells.Item(24, 3).Value = 10.05519576933455
$ws.Cells.Item(24, 4).Value = 5.134861229890784
$ws.Cells.Item(24, 6).Value = 22.87533506448239
$ws.Cells.Item(24, 7).Value = 26.27354267716675
$ws.Cells.Item(24, 8).Value = 13.64826080461569
$ws.Cells.Item(24, 9).Value = 21.51312616981474
$ws.Cells.Item(24, 11).Value = 8.812554224363442
$ws.Cells.Item(24, 12).Value = 10.80183916405954
$ws.Cells.Item(24, 14).Value = 18.03430458365987
$ws.Cells.Item(24, 15).Value = 20.44568270703185

$ws.Cells.Item(25, 2).Value = 11.48531664941857
$ws.Cells.Item(25, 3).Value = 10.0869370068084
$ws.Cells.Item(25, 4).Value = 4.983699769639493
$ws.Cells.Item(25, 6).Value = 22.88763945284946
$ws.Cells.Item(25, 7).Value = 26.30991695576829
$ws.Cells.Item(25, 8).Value = 13.71201554693133
$ws.Cells.Item(25, 9).Value = 21.6377825116794
$ws.Cells.Item(25, 11).Value = 8.315866667753971
$ws.Cells.Item(25, 12).Value = 10.72700048965314
$ws.Cells.Item(25, 14).Value = 18.12362326732458
$ws.Cells.Item(25, 15).Value = 20.53802750922556
